$d = $word.ActiveDocument

$replacements = @(
    @{old="71×19=1349"; new="79×59=4661"},
    @{old="75×74=5550"; new="25×91=2275"},
    @{old="69×52=3588"; new="88×93=8184"},
    @{old="92×40=3680"; new="49×64=3136"},
    @{old="35×24=840"; new="86×15=1290"},
    @{old="92×67=6164"; new="26×58=1508"},
    @{old="97×31=3007"; new="30×28=840"},
    @{old="96×84=8064"; new="95×63=5985"},
    @{old="89×20=1780"; new="40×24=960"},
    @{old="50×54=2700"; new="45×98=4410"},
    @{old="67×12=804"; new="57×32=1824"},
    @{old="83×85=7055"; new="79×34=2686"},
    @{old="22×17=374"; new="33×79=2607"},
    @{old="25×58=1450"; new="52×25=1300"},
    @{old="68×56=3808"; new="30×59=1770"},
    @{old="64×39=2496"; new="59×74=4366"},
    @{old="30×76=2280"; new="65×76=4940"},
    @{old="15×38=570"; new="83×64=5312"},
    @{old="25×71=1775"; new="67×15=1005"},
    @{old="47×84=3948"; new="55×28=1540"},
    @{old="44×71=3124"; new="36×87=3132"},
    @{old="35×88=3080"; new="93×38=3534"},
    @{old="24×61=1464"; new="13×16=208"},
    @{old="30×54=1620"; new="20×39=780"},
    @{old="67×32=2144"; new="72×18=1296"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
